# Apply "Add Leave Card 10/32023 3:18 PM" edits
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws2 = $wb.Worksheets.Item("2017 LEAVE BALANCE")

# ---- Sheet 1: "2018 LEAVE CREDITS" (Table13) ----

# Row 87 - EARNED 1.25 (period already 45047 / PARTICULARS already set)
$ws1.Range("C87").Value = 1.25

# Row 89
$ws1.Range("A89").Value = 45078
$ws1.Range("C89").Value = 1.25

# Row 90
$ws1.Range("A90").Value = 45108
$ws1.Range("C90").Value = 1.25

# Row 91
$ws1.Range("A91").Value = 45139
$ws1.Range("C91").Value = 1.25

# Row 92
$ws1.Range("A92").Value = 45170
$ws1.Range("B92").Value = "SL(2-0-0)"
$ws1.Range("C92").Value = 1.25
$ws1.Range("H92").Value = 2
$ws1.Range("K92").Value = "9/11,12/2023"

# Row 93
$ws1.Range("B93").Value = "SL(1-0-0)"
$ws1.Range("H93").Value = 1
$ws1.Range("K93").Value = 45191

# Update BALANCE formulas on row 9 (drop the CONVERTION offset)
$ws1.Range("E9").Formula = "=SUM(Table13[EARNED])-SUM(Table13[Absence Undertime W/ Pay])"
$ws1.Range("I9").Formula = "=SUM(Table13[[EARNED ]])-SUM(Table13[Absence Undertime  W/ Pay])"

# ---- Sheet 2: "2017 LEAVE BALANCE" (Table1) ----

# Row 47
$ws2.Range("A47").Value = 45139
$ws2.Range("B47").Value = "VL(1-0-0)"
$ws2.Range("D47").Value = 1
$ws2.Range("K47").Value = 45160

# Row 48
$ws2.Range("A48").Value = 45170
$ws2.Range("B48").Value = "VL(1-0-0)"
$ws2.Range("D48").Value = 1
$ws2.Range("K48").Value = 45190

# ---- View/selection state ----
$ws1.Range("I10").Select()
$ws2.Range("B49").Select()

$ws2.Activate()
$ws2.Range("B3:C3").Select()
